# DPLKINV118-001..003 - Jenis Porto - Reksadana Approve Data
# Scripting DPLKKPS041-001 until DPLKKPS096-001 (14/02/2023)
#
# Change "Kode Jenis Porto" from R07 to R91 on each of the three test-case
# sheets, and leave the workbook with the last sheet (DPLKINV118-003)
# active/selected, matching the cell selections recorded in the edited file.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("DPLKINV118-001")
$ws2 = $wb.Worksheets.Item("DPLKINV118-002")
$ws3 = $wb.Worksheets.Item("DPLKINV118-003")

# KODE_JENIS_PORTO column (M) : R07 -> R91
$ws1.Range("M2").Value = "R91"
$ws2.Range("M2").Value = "R91"
$ws3.Range("M2").Value = "R91"

# Restore each sheet's on-screen selection
$ws1.Activate()
$ws1.Range("N2").Select() | Out-Null

$ws2.Activate()
$ws2.Range("N2").Select() | Out-Null

$ws3.Activate()
$ws3.Range("O2").Select() | Out-Null
